$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.554428100585938
$ws.Range("B1").Value = 2.831684589385986
$ws.Range("C1").Value = 2.110269546508789
$ws.Range("D1").Value = 1.888778686523438
$ws.Range("E1").Value = 1.661590337753296
